# Auto-generated script applying cell value updates described in the commit diff
# for Sheets/Behemoth_Profits.xlsx (source workbook with 8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39499.5
$ws.Range("J3").Value = 39499.5
$ws.Range("L3").Value = 39499.5
$ws.Range("N3").Value = -39727.5
$ws.Range("H9").Value = 337.75
$ws.Range("I9").Value = 314.77777
$ws.Range("K9").Value = 314.77777
$ws.Range("M9").Value = -145.77777
$ws.Range("H29").Value = 1149.6666
$ws.Range("J29").Value = 1224.5
$ws.Range("L29").Value = 3673.5
$ws.Range("N29").Value = -4235.5
$ws.Range("H31").Value = 2551
$ws.Range("J31").Value = 4999
$ws.Range("L31").Value = 14997
$ws.Range("N31").Value = -15457
$ws.Range("H33").Value = 2094.2917
$ws.Range("I33").Value = 1758.4736
$ws.Range("K33").Value = 1758.4736
$ws.Range("M33").Value = -1529.4736
$ws.Range("H46").Value = 2200
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2200
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 6600
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6838
$ws.Range("H60").Value = 2200
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 2200
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 6600
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -7568
$ws.Range("H62").Value = 14227.667
$ws.Range("I62").Value = 4920
$ws.Range("J62").Value = 25862.25
$ws.Range("K62").Value = 4920
$ws.Range("L62").Value = 25862.25
$ws.Range("M62").Value = -4296
$ws.Range("N62").Value = -27110.25
$ws.Range("H65").Value = 14227.667
$ws.Range("I65").Value = 4920
$ws.Range("J65").Value = 25862.25
$ws.Range("K65").Value = 24600
$ws.Range("L65").Value = 129311.25
$ws.Range("M65").Value = -21480
$ws.Range("N65").Value = -135551.25
$ws.Range("H74").Value = 2249.25
$ws.Range("I74").Value = 2332.3333
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 2332.3333
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -1396.3333
$ws.Range("N74").Value = -3872
$ws.Range("H77").Value = 2249.25
$ws.Range("I77").Value = 2332.3333
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 11661.6665
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -6981.666499999999
$ws.Range("N77").Value = -19360
$ws.Range("H102").Value = 39499.5
$ws.Range("J102").Value = 39499.5
$ws.Range("L102").Value = 39499.5
$ws.Range("N102").Value = -45989.5
$ws.Range("H106").Value = 5862.375
$ws.Range("I106").Value = 5452.5
$ws.Range("K106").Value = 5452.5
$ws.Range("M106").Value = -4821.5
$ws.Range("H112").Value = 1625.12
$ws.Range("J112").Value = 1866.5555
$ws.Range("L112").Value = 5599.666499999999
$ws.Range("N112").Value = -7815.666499999999
$ws.Range("H113").Value = 38464864
$ws.Range("I113").Value = 12502722
$ws.Range("K113").Value = 12502722
$ws.Range("M113").Value = -12499468
$ws.Range("H125").Value = 2299
$ws.Range("I125").Value = 1266.3334
$ws.Range("K125").Value = 11397.0006
$ws.Range("M125").Value = -8937.000599999999
$ws.Range("H132").Value = 1470.16
$ws.Range("I132").Value = 1470.16
$ws.Range("K132").Value = 4410.48
$ws.Range("M132").Value = -1880.48
$ws.Range("H137").Value = 11187.462
$ws.Range("J137").Value = 22166.5
$ws.Range("L137").Value = 66499.5
$ws.Range("N137").Value = -71599.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 216
$ws.Range("I22").Value = 216
$ws.Range("K22").Value = 216
$ws.Range("M22").Value = 83
$ws.Range("H39").Value = 4333.3335
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H45").Value = 27780194
$ws.Range("I45").Value = 45456080
$ws.Range("J45").Value = 3799.8572
$ws.Range("K45").Value = 45456080
$ws.Range("L45").Value = 3799.8572
$ws.Range("M45").Value = -45455703
$ws.Range("N45").Value = -4553.8572
$ws.Range("H61").Value = 30008798
$ws.Range("I61").Value = 23817594
$ws.Range("K61").Value = 23817594
$ws.Range("M61").Value = -23817382
$ws.Range("H75").Value = 45000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 45000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 45000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -46748
$ws.Range("H78").Value = 45000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 45000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 135000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -143736
$ws.Range("H97").Value = 2025.3478
$ws.Range("I97").Value = 1994.75
$ws.Range("J97").Value = 2041.6666
$ws.Range("K97").Value = 1994.75
$ws.Range("L97").Value = 2041.6666
$ws.Range("M97").Value = -1498.75
$ws.Range("N97").Value = -3033.6666
$ws.Range("H107").Value = 109000
$ws.Range("J107").Value = 109000
$ws.Range("L107").Value = 109000
$ws.Range("N107").Value = -116680
$ws.Range("H122").Value = 2126.75
$ws.Range("I122").Value = 1993
$ws.Range("J122").Value = 2171.3333
$ws.Range("K122").Value = 5979
$ws.Range("L122").Value = 6513.999899999999
$ws.Range("M122").Value = -3529
$ws.Range("N122").Value = -11413.9999
$ws.Range("H132").Value = 3691.2927
$ws.Range("I132").Value = 1963.9286
$ws.Range("J132").Value = 7411.769
$ws.Range("K132").Value = 5891.7858
$ws.Range("L132").Value = 22235.307
$ws.Range("M132").Value = -3361.7858
$ws.Range("N132").Value = -27295.307
$ws.Range("H136").Value = 30008798
$ws.Range("I136").Value = 23817594
$ws.Range("K136").Value = 71452782
$ws.Range("M136").Value = -71450232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 3349.8
$ws.Range("I25").Value = 4024.75
$ws.Range("J25").Value = 650
$ws.Range("K25").Value = 4024.75
$ws.Range("L25").Value = 650
$ws.Range("M25").Value = -3789.75
$ws.Range("N25").Value = -1120
$ws.Range("H37").Value = 6428.625
$ws.Range("I37").Value = 1880
$ws.Range("J37").Value = 14009.667
$ws.Range("K37").Value = 1880
$ws.Range("L37").Value = 14009.667
$ws.Range("M37").Value = -1743
$ws.Range("N37").Value = -14283.667
$ws.Range("H46").Value = 4749.5
$ws.Range("J46").Value = 4749.5
$ws.Range("L46").Value = 4749.5
$ws.Range("N46").Value = -5345.5
$ws.Range("H99").Value = 3902.9412
$ws.Range("I99").Value = 3525
$ws.Range("J99").Value = 4109.091
$ws.Range("K99").Value = 3525
$ws.Range("L99").Value = 4109.091
$ws.Range("M99").Value = -2027
$ws.Range("N99").Value = -7105.091
$ws.Range("H112").Value = 99995
$ws.Range("J112").Value = 99995
$ws.Range("L112").Value = 99995
$ws.Range("N112").Value = -102949

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1566.9
$ws.Range("I16").Value = 1460.8334
$ws.Range("K16").Value = 1460.8334
$ws.Range("M16").Value = -1173.8334
$ws.Range("H20").Value = 120000
$ws.Range("J20").Value = 120000
$ws.Range("L20").Value = 120000
$ws.Range("N20").Value = -120472
$ws.Range("H30").Value = 120000
$ws.Range("J30").Value = 120000
$ws.Range("L30").Value = 120000
$ws.Range("N30").Value = -120182
$ws.Range("H31").Value = 606496.4
$ws.Range("J31").Value = 711821.5
$ws.Range("L31").Value = 711821.5
$ws.Range("N31").Value = -712411.5
$ws.Range("H34").Value = 606496.4
$ws.Range("J34").Value = 711821.5
$ws.Range("L34").Value = 711821.5
$ws.Range("N34").Value = -712225.5
$ws.Range("H38").Value = 500
$ws.Range("I38").Value = 500
$ws.Range("K38").Value = 500
$ws.Range("M38").Value = -123
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("K46").Value = 500
$ws.Range("M46").Value = -289
$ws.Range("H58").Value = 3113.75
$ws.Range("I58").Value = 1047.25
$ws.Range("J58").Value = 7246.75
$ws.Range("K58").Value = 1047.25
$ws.Range("L58").Value = 7246.75
$ws.Range("M58").Value = -844.25
$ws.Range("N58").Value = -7652.75
$ws.Range("H86").Value = 8332.333000000001
$ws.Range("I86").Value = 8499.5
$ws.Range("K86").Value = 8499.5
$ws.Range("M86").Value = -7376.5
$ws.Range("H89").Value = 8332.333000000001
$ws.Range("I89").Value = 8499.5
$ws.Range("K89").Value = 42497.5
$ws.Range("M89").Value = -36881.5
$ws.Range("H99").Value = 2199.125
$ws.Range("I99").Value = 1701.9166
$ws.Range("K99").Value = 1701.9166
$ws.Range("M99").Value = -203.9166
$ws.Range("H105").Value = 1261.2778
$ws.Range("I105").Value = 1294.2
$ws.Range("K105").Value = 1294.2
$ws.Range("M105").Value = 452.8
$ws.Range("H107").Value = 2824.2
$ws.Range("I107").Value = 1011
$ws.Range("J107").Value = 3277.5
$ws.Range("K107").Value = 1011
$ws.Range("L107").Value = 3277.5
$ws.Range("M107").Value = 909
$ws.Range("N107").Value = -7117.5
$ws.Range("H113").Value = 1566.9
$ws.Range("I113").Value = 1460.8334
$ws.Range("K113").Value = 1460.8334
$ws.Range("M113").Value = 709.1666
$ws.Range("H122").Value = 1941
$ws.Range("J122").Value = 1911
$ws.Range("L122").Value = 5733
$ws.Range("N122").Value = -10633
$ws.Range("H126").Value = 2199.125
$ws.Range("I126").Value = 1701.9166
$ws.Range("K126").Value = 5105.7498
$ws.Range("M126").Value = -2635.7498
$ws.Range("H128").Value = 120000
$ws.Range("J128").Value = 120000
$ws.Range("L128").Value = 120000
$ws.Range("N128").Value = -129960
$ws.Range("H136").Value = 3113.75
$ws.Range("I136").Value = 1047.25
$ws.Range("J136").Value = 7246.75
$ws.Range("K136").Value = 3141.75
$ws.Range("L136").Value = 21740.25
$ws.Range("M136").Value = -591.75
$ws.Range("N136").Value = -26840.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10091108
$ws.Range("I4").Value = 7562637
$ws.Range("J4").Value = 16833696
$ws.Range("K4").Value = 22687911
$ws.Range("L4").Value = 50501088
$ws.Range("M4").Value = -22687799
$ws.Range("N4").Value = -50501312
$ws.Range("H7").Value = 391.91666
$ws.Range("I7").Value = 444.83334
$ws.Range("J7").Value = 339
$ws.Range("K7").Value = 1334.50002
$ws.Range("L7").Value = 1017
$ws.Range("M7").Value = -1222.50002
$ws.Range("N7").Value = -1241
$ws.Range("H31").Value = 971.75
$ws.Range("I31").Value = 950
$ws.Range("J31").Value = 993.5
$ws.Range("K31").Value = 2850
$ws.Range("L31").Value = 2980.5
$ws.Range("M31").Value = -2562
$ws.Range("N31").Value = -3556.5
$ws.Range("H46").Value = 1399.8
$ws.Range("I46").Value = 333
$ws.Range("K46").Value = 999
$ws.Range("M46").Value = -908
$ws.Range("H50").Value = 425
$ws.Range("I50").Value = 268
$ws.Range("K50").Value = 804
$ws.Range("M50").Value = -323
$ws.Range("H53").Value = 425
$ws.Range("I53").Value = 268
$ws.Range("K53").Value = 804
$ws.Range("M53").Value = -323
$ws.Range("H97").Value = 8930147
$ws.Range("I97").Value = 17858192
$ws.Range("J97").Value = 2102
$ws.Range("K97").Value = 53574576
$ws.Range("L97").Value = 6306
$ws.Range("M97").Value = -53574080
$ws.Range("N97").Value = -7298
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 133.5238
$ws.Range("I2").Value = 39
$ws.Range("J2").Value = 237.5
$ws.Range("K2").Value = 39
$ws.Range("L2").Value = 237.5
$ws.Range("M2").Value = 74
$ws.Range("N2").Value = -463.5
$ws.Range("H31").Value = 2770
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 2770
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H41").Value = 12500
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15710
$ws.Range("H102").Value = 1310.1818
$ws.Range("I102").Value = 1310.1818
$ws.Range("K102").Value = 1310.1818
$ws.Range("M102").Value = 311.8181999999999
$ws.Range("H122").Value = 2950
$ws.Range("I122").Value = 2933.3333
$ws.Range("K122").Value = 8799.999899999999
$ws.Range("M122").Value = -6349.999899999999
$ws.Range("H126").Value = 6666
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 20000
$ws.Range("J25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("N25").Value = -20460
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H122").Value = 5521.6665
$ws.Range("I122").Value = 4813.5713
$ws.Range("K122").Value = 14440.7139
$ws.Range("M122").Value = -11990.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 129900
$ws.Range("J4").Value = 1500
$ws.Range("L4").Value = 1500
$ws.Range("N4").Value = -1726
$ws.Range("H39").Value = 27013.334
$ws.Range("J39").Value = 27747.5
$ws.Range("L39").Value = 27747.5
$ws.Range("N39").Value = -28573.5
$ws.Range("H42").Value = 58989
$ws.Range("I42").Value = 58989
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 58989
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -58611
$ws.Range("N42").ClearContents()
$ws.Range("H96").Value = 4282.357
$ws.Range("J96").Value = 4317.6665
$ws.Range("L96").Value = 4317.6665
$ws.Range("N96").Value = -7063.6665
$ws.Range("H122").Value = 1408.9412
$ws.Range("I122").Value = 1397
$ws.Range("K122").Value = 4191
$ws.Range("M122").Value = -1741
$ws.Range("H126").Value = 2830.818
$ws.Range("I126").Value = 3057.6667
$ws.Range("K126").Value = 9173.000100000001
$ws.Range("M126").Value = -6703.000100000001
$ws.Range("H135").Value = 99178.5
$ws.Range("J135").Value = 99178.5
$ws.Range("L135").Value = 99178.5
$ws.Range("N135").Value = -109318.5
